# Update the LR-pairs worksheet with new TPM-based values.
# Adds a "MuSCs" sending-cluster block (rows 7-11) mirroring the existing
# "ECs" sending-cluster block (rows 2-6), and refreshes all numeric columns
# (E:T) for every row with newly-computed TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    # Row, A(Sending cluster), B(Ligand), C(Receptor), D(Target cluster), E..T
    @(2,  "ECs",   "Fgf9", "Fgfr1", "ECs",              2, 0.6666666666666666, 0.7364236666666667, 2.209271,           0.9917500467982164, 0.9944848646626661, 3, 1, 5.978421000000001,  17.935263,           0.05704457007880161,  0.06242884486533885,  4.402650713697001,   39.62385642327301,  0.05657395504523563,  0.06208454133695308),
    @(3,  "ECs",   "Fgf9", "Fgfr1", "FAPs",             2, 0.6666666666666666, 0.7364236666666667, 2.209271,           0.9917500467982164, 0.9944848646626661, 3, 1, 69.657865,           208.973595,          0.6646576013185088,   0.7273927426214574,   51.29770035547167,   461.679303199245,   0.6591742072124214,   0.7233810732025056),
    @(4,  "ECs",   "Fgf9", "Fgfr1", "Inflammatory-Mac", 2, 0.6666666666666666, 0.7364236666666667, 2.209271,           0.9917500467982164, 0.9944848646626661, 3, 1, 1.290243,            3.870729,            0.01231116999491725,  0.01347318632889677,  0.9501654809510001,  8.551489328559001,  0.01220960341859999,  0.01339887988286779),
    @(5,  "ECs",   "Fgf9", "Fgfr1", "MuSCs",            2, 0.6666666666666666, 0.7364236666666667, 2.209271,           0.9917500467982164, 0.9944848646626661, 2, 1, 27.1166075,          54.233215,           0.2587397603536297,   0.1887743138075849,   19.96931152271084,   119.815869136265,   0.2566051694392715,   0.1877331979187238),
    @(6,  "ECs",   "Fgf9", "Fgfr1", "Resolving-Mac",    2, 0.6666666666666666, 0.7364236666666667, 2.209271,           0.9917500467982164, 0.9944848646626661, 3, 1, 0.759494,            2.278482,            0.00724689825414258,  0.007930912376722157, 0.5593093562913334,  5.033784206622,     0.007187111682687817, 0.007887172321615997),
    @(7,  "MuSCs", "Fgf9", "Fgfr1", "ECs",              1, 0.5,                0.006126,            0.012252,           0.008249953201783585, 0.005515135337333892, 3, 1, 5.978421000000001,  17.935263,           0.05704457007880161,  0.06242884486533885,  0.03662380704600001, 0.219742842276,     0.0004706150335659774, 0.0003443035283857658),
    @(8,  "MuSCs", "Fgf9", "Fgfr1", "FAPs",             1, 0.5,                0.006126,            0.012252,           0.008249953201783585, 0.005515135337333892, 3, 1, 69.657865,           208.973595,          0.6646576013185088,   0.7273927426214574,   0.42672408099,       2.56034448594,      0.00548339410608743,  0.004011669418951816),
    @(9,  "MuSCs", "Fgf9", "Fgfr1", "Inflammatory-Mac", 1, 0.5,                0.006126,            0.012252,           0.008249953201783585, 0.005515135337333892, 3, 1, 1.290243,            3.870729,            0.01231116999491725,  0.01347318632889677,  0.007904028618000001, 0.047424171708,    0.0001015665763172696, 0.0000743064460289825),
    @(10, "MuSCs", "Fgf9", "Fgfr1", "MuSCs",            1, 0.5,                0.006126,            0.012252,           0.008249953201783585, 0.005515135337333892, 2, 1, 27.1166075,          54.233215,           0.2587397603536297,   0.1887743138075849,   0.166116337545,      0.6644653501800001, 0.002134590914358145, 0.001041115888861169),
    @(11, "MuSCs", "Fgf9", "Fgfr1", "Resolving-Mac",    1, 0.5,                0.006126,            0.012252,           0.008249953201783585, 0.005515135337333892, 3, 1, 0.759494,            2.278482,            0.00724689825414258,  0.007930912376722157, 0.004652660244,      0.027915961464,     0.00005978657145476345, 0.00004374005510615909)
)

foreach ($row in $data) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]

    for ($c = 5; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
}
